$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + week dates), preserving rich-text runs ---
$a8 = $ws.Range("A8")
$a8v = $a8.Value2
$a8.Characters($a8v.Length, 1).Text = "3"
$a8c = $a8.Characters($a8v.Length, 1)
$a8c.Font.Name = "Andale WT"
$a8c.Font.Size = 10

$c9 = $ws.Range("C9")
# "Report Covering the Week  1/6/2025  Through  1/12/2025"
# run2 = "1/6/2025" (8 chars) starting right after "Report Covering the Week  " (27 chars)
$run2Start = 28
$run2Len = 8
$c9.Characters($run2Start, $run2Len).Text = "1/13/2025"
$run2b = $c9.Characters($run2Start, 9)
$run2b.Font.Name = "Andale WT"
$run2b.Font.Size = 10

$v2 = $c9.Value2
# run4 = "1/12/2025" is the last 10 characters
$run4Len = 10
$run4Start = $v2.Length - $run4Len + 1
$c9.Characters($run4Start, $run4Len).Text = "1/19/2025"
$run4b = $c9.Characters($run4Start, 9)
$run4b.Font.Name = "Andale WT"
$run4b.Font.Size = 10

# --- Column E width (auto bestFit changed because of wider percentage value in E17) ---
$ws.Columns.Item(5).ColumnWidth = 8.43

# --- Cells changing from N/A placeholder text to numeric values ---
$donor14 = $ws.Range("F15")
$donor15 = $ws.Range("H14")
$donor14.Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Value = 1
$donor15.Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E14").Value = -100
$donor14.Copy()
$ws.Range("J14").PasteSpecial(-4122)
$ws.Range("J14").Value = 1
$donor15.Copy()
$ws.Range("K14").PasteSpecial(-4122)
$ws.Range("K14").Value = -100
$donor15.Copy()
$ws.Range("L15").PasteSpecial(-4122)
$ws.Range("L15").Value = 0
$donor15.Copy()
$ws.Range("N15").PasteSpecial(-4122)
$ws.Range("N15").Value = -50
$donor14.Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$donor15.Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$donor14.Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("G22").Value = 1
$donor15.Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("H22").Value = 0
$donor14.Copy()
$ws.Range("J22").PasteSpecial(-4122)
$ws.Range("J22").Value = 1
$donor15.Copy()
$ws.Range("K22").PasteSpecial(-4122)
$ws.Range("K22").Value = -100
$donor15.Copy()
$ws.Range("L27").PasteSpecial(-4122)
$ws.Range("L27").Value = 0
$donor15.Copy()
$ws.Range("N29").PasteSpecial(-4122)
$ws.Range("N29").Value = -100
$donor15.Copy()
$ws.Range("N30").PasteSpecial(-4122)
$ws.Range("N30").Value = -100
$excel.CutCopyMode = $false

# --- Cells changing from numeric values to N/A placeholder text ---
$donor13 = $ws.Range("A14")
$donor13.Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Formula = "=""0"""
$ws.Range("C15").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$donor13.Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Formula = "=""0"""
$ws.Range("C27").Copy()
$ws.Range("C27").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("G14").Value = 2
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 75
$ws.Range("F16").Value = 20
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 16
$ws.Range("J16").Value = 5
$ws.Range("K16").Value = 220
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = -20
$ws.Range("N16").Value = -67.346938775510
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 26
$ws.Range("H17").Value = 3.846153846153
$ws.Range("I17").Value = 16
$ws.Range("J17").Value = 10
$ws.Range("K17").Value = 60
$ws.Range("L17").Value = -27.272727272727
$ws.Range("M17").Value = 166.666666666667
$ws.Range("N17").Value = 23.076923076923
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -38.888888888888
$ws.Range("I18").Value = 7
$ws.Range("J18").Value = 13
$ws.Range("K18").Value = -46.153846153846
$ws.Range("M18").Value = -46.153846153846
$ws.Range("N18").Value = -90.789473684210
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = 3.125
$ws.Range("I19").Value = 24
$ws.Range("J19").Value = 18
$ws.Range("K19").Value = 33.333333333333
$ws.Range("L19").Value = -31.428571428571
$ws.Range("M19").Value = 71.428571428571
$ws.Range("N19").Value = -25
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -60
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 33
$ws.Range("H20").Value = -39.393939393939
$ws.Range("I20").Value = 11
$ws.Range("J20").Value = 26
$ws.Range("K20").Value = -57.692307692307
$ws.Range("L20").Value = -50
$ws.Range("M20").Value = -38.888888888888
$ws.Range("N20").Value = -93.529411764705
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 0
$ws.Range("G21").Value = 121
$ws.Range("H21").Value = -6.611570247933
$ws.Range("I21").Value = 76
$ws.Range("J21").Value = 73
$ws.Range("K21").Value = 4.109589041095
$ws.Range("L21").Value = -20.833333333333
$ws.Range("M21").Value = 7.042253521126
$ws.Range("N21").Value = -78.097982708933
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = 138.461538461538
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 77
$ws.Range("H24").Value = 24.675324675324
$ws.Range("I24").Value = 65
$ws.Range("J24").Value = 51
$ws.Range("K24").Value = 27.450980392156
$ws.Range("L24").Value = 20.370370370370
$ws.Range("M24").Value = 103.125
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 220
$ws.Range("F25").Value = 50
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = 92.307692307692
$ws.Range("I25").Value = 34
$ws.Range("J25").Value = 15
$ws.Range("K25").Value = 126.666666666667
$ws.Range("L25").Value = 30.769230769230
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 52
$ws.Range("G26").Value = 39
$ws.Range("H26").Value = 33.333333333333
$ws.Range("I26").Value = 31
$ws.Range("J26").Value = 27
$ws.Range("K26").Value = 14.814814814814
$ws.Range("L26").Value = 24
$ws.Range("M26").Value = 55
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 50
$ws.Range("J44").Value = 526
$ws.Range("K44").Value = 38.057742782152
$ws.Range("L44").Value = 19.545454545454
$ws.Range("M44").Value = -11.596638655462
$ws.Range("N44").Value = -26.433566433566
$ws.Range("J46").Value = 1609
$ws.Range("K46").Value = -26.327838827838
$ws.Range("L46").Value = -52.127343052662
$ws.Range("M46").Value = -75.177414378278
$ws.Range("N46").Value = -74.625453398517
